$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 7: Mecânica - Campus Congonhas
$ws.Range("E7").Value = 18
$ws.Range("F7").Value = 11
$ws.Range("H7").Value = 11

# Row 8: Mineração - Campus Congonhas
$ws.Range("F8").Value = 11
$ws.Range("H8").Value = 11

# Row 9: Eletrotécnica - Campus Conselheiro Lafaiete
$ws.Range("E9").Value = 13

# Row 15: Metalurgia - Campus Ouro Preto
$ws.Range("E15").Value = 63

# Row 16: Mineração - Campus Ouro Preto
$ws.Range("F16").Value = 63
$ws.Range("H16").Value = 63
